# Auto update Excel log 2026-02-04 14:07:46
# Appends new sensor-log rows to the PIR, Humidity and Temperature sheets.

$wb = $excel.ActiveWorkbook

function Append-LogRows($ws, $startRow, $rows) {
    $endRow = $startRow + $rows.Count - 1
    $rng = $ws.Range("A$startRow" + ":F$endRow")
    # Force text storage so values like "2026-02-04", "14:00" and "76.6%"
    # are kept as literal strings instead of being auto-converted to
    # dates/times/percentages by Excel's smart-entry heuristics.
    $rng.NumberFormat = "@"

    $r = $startRow
    foreach ($row in $rows) {
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $ws.Cells.Item($r, 5).Value = $row[4]
        $ws.Cells.Item($r, 6).Value = $row[5]
        $r++
    }

    # Drop the explicit "Text" number format again so the new cells end up
    # using the workbook's default (unstyled) cell format, matching the
    # rest of the log.
    $rng.ClearFormats()
}

# ---------------------------------------------------------------------------
# PIR sheet - append rows 67-79
# ---------------------------------------------------------------------------
$wsPIR = $wb.Worksheets.Item("PIR")
$pirRows = @(
    @("2026-02-04","14:06:29","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:06:30","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:06:30","14:00","Bathroom","Motion Detected","Active"),
    @("2026-02-04","14:06:38","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:06:39","14:00","Bathroom","Motion Detected","Active"),
    @("2026-02-04","14:06:46","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:06:49","14:00","Bathroom","Motion Detected","Active"),
    @("2026-02-04","14:06:57","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:07:02","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:07:07","14:00","Bathroom","Motion Detected","Active"),
    @("2026-02-04","14:07:15","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:07:19","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:07:24","14:00","Bathroom","No Motion","Inactive")
)
Append-LogRows $wsPIR 67 $pirRows

# ---------------------------------------------------------------------------
# Humidity sheet - append rows 56-62
# ---------------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @("2026-02-04","14:06:28","14:00","Bathroom","76.6%","Active"),
    @("2026-02-04","14:06:34","14:00","Bathroom","76.4%","Active"),
    @("2026-02-04","14:06:49","14:00","Bathroom","76.5%","Active"),
    @("2026-02-04","14:06:59","14:00","Bathroom","76.5%","Active"),
    @("2026-02-04","14:07:09","14:00","Bathroom","77.6%","Active"),
    @("2026-02-04","14:07:19","14:00","Bathroom","77.7%","Active"),
    @("2026-02-04","14:07:24","14:00","Bathroom","77.6%","Active")
)
Append-LogRows $wsHumidity 56 $humidityRows

# ---------------------------------------------------------------------------
# Temperature sheet - append rows 56-62
# ---------------------------------------------------------------------------
$wsTemperature = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    @("2026-02-04","14:06:28","14:00","Bathroom","24.8C","Active"),
    @("2026-02-04","14:06:34","14:00","Bathroom","24.8C","Active"),
    @("2026-02-04","14:06:49","14:00","Bathroom","24.8C","Active"),
    @("2026-02-04","14:06:59","14:00","Bathroom","24.8C","Active"),
    @("2026-02-04","14:07:09","14:00","Bathroom","24.8C","Active"),
    @("2026-02-04","14:07:19","14:00","Bathroom","24.8C","Active"),
    @("2026-02-04","14:07:24","14:00","Bathroom","24.8C","Active")
)
Append-LogRows $wsTemperature 56 $temperatureRows
